$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "Texas Notes" sheet by duplicating "Calculations" (keeps
#    number formats / fonts / column widths identical to the source sheet)
#    and position it right after "Calculations".
# ---------------------------------------------------------------------------
$calc = $wb.Worksheets.Item("Calculations")
$calc.Copy([System.Reflection.Missing]::Value, $calc)
$tx = $wb.Worksheets.Item("Calculations (2)")
$tx.Name = "Texas Notes"

# Push the existing table down 4 rows so it lands on rows 5-10 (as in the
# target layout) and leaves room for the two intro note lines at the top.
$tx.Rows("1:4").Insert()

# Row 1-2: intro notes about the Texas-specific graphical data.
$tx.Range("A1").Value = "The source has Texas specific data, but it's in graphical form. "
$tx.Range("A2").Value = "So, I used some visual editing software to measure the areas of the Texas specific data. This has some error associated with it."

# Row 5 headers: relabel the BAU / Extended-ITC columns as pixel measurements.
$tx.Range("B5").Value = "BAU Deployment 2015-2022 (square pixels measured)"
$tx.Range("C5").Value = "Deployment with Extended ITC 2015-2022 (square pixels measured)"

# Rows 6-7: Texas-specific pixel-area data (replacing the national GW figures).
$tx.Range("B6").Value = 26477
$tx.Range("C6").Value = 29137
$tx.Range("B7").Value = 11010
$tx.Range("C7").Value = 20634

# Rows 16-17: closing notes about averaging with the national numbers.
$tx.Range("A16").Value = "The main point here is that Texas residential is less elastic and commericial is more elastic than the national average. Because the measuring technique I used has some error in it, I will average the numbres"
$tx.Range("A17").Value = "above with the national numbers in the ""Calculations"" tab to come up with something a bit more conservative, in case my measuring error is high."

# Rows 19-20: average the Texas-measured elasticities with the national ones.
$tx.Range("A19").Value = "Residential"
$tx.Range("B19").Formula = "=AVERAGE(B13,Calculations!B9)"
$tx.Range("A20").Value = "Commercial"
$tx.Range("B20").Formula = "=AVERAGE(B14,Calculations!B10)"

$tx.Range("B21").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Point "EoDSDwSP" at the new "Texas Notes" tab instead of "Calculations".
# ---------------------------------------------------------------------------
$eods = $wb.Worksheets.Item("EoDSDwSP")
$eods.Range("B2").Formula = "='Texas Notes'!B19"
$eods.Range("B4").Formula = "='Texas Notes'!B20"
$eods.Range("B5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Add the hyperlink on the "About" sheet's source-url cell.
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")
$urlCell = $about.Range("B6")
$about.Hyperlinks.Add($urlCell, $urlCell.Text) | Out-Null
$urlCell.Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Restore view selections on the other sheets, then leave EoDSDwSP active.
# ---------------------------------------------------------------------------
$calc.Range("B9").Select() | Out-Null

$eods.Activate()
